$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.665.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.233.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'610.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.15%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'159.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.231.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.84%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.72%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.71%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.99%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'38.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.32%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.763.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'66.720.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.33%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.238.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.98%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.39%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'511.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.48%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.28%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'14.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.30%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'85.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.15%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'3.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.13%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.90%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +2.83%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.39%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +33.95%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'28.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.08%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.99%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.17%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'508.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.71%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'55.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.35%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +16.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +7.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +6.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.36%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.06%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.906.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.87%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'28.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.72%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +4.30%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = "'  -0.39%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'123.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.07%  "
$ws.Range("E51").Style = "Normal"
